$d = $word.ActiveDocument

# 1) Update image relationship target rId5 -> rId7 is handled internally by Word when
#    parts are renumbered; we don't set this directly.

# 2) Remove the _GoBack bookmark from the body paragraph (it is being moved to the footer)
foreach ($bm in $d.Bookmarks) {
    if ($bm.Name -eq "_GoBack") {
        $bm.Delete()
    }
}

# 3) Set up different headers/footers for even/odd/first pages
$sec = $d.Sections(1)

# Unlink footers from "previous" (there's no previous, but this forces per-section content)
$ftrPrimary = $sec.Footers(1)
$ftrFirst = $sec.Footers(2)
$ftrEven = $sec.Footers(3)

$hdrPrimary = $sec.Headers(1)
$hdrFirst = $sec.Headers(2)
$hdrEven = $sec.Headers(3)

$ftrPrimary.LinkToPrevious = $false
$ftrFirst.LinkToPrevious = $false
$ftrEven.LinkToPrevious = $false

$hdrPrimary.LinkToPrevious = $false
$hdrFirst.LinkToPrevious = $false
$hdrEven.LinkToPrevious = $false

# Move the bookmark into the default footer
$ftrPrimary.Range.Text = ""
$d.Bookmarks.Add("_GoBack", $ftrPrimary.Range)
